$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Notes sheet: select the existing last history row first (B27 is picked
# later after the new row is appended, but we record the pre-edit selection
# state on this sheet so its own tabSelected flag eventually gets dropped
# once sheet1 is reactivated at the end of the script).

# Copy formatting of the last history row (row 25) down into the new row 26
# so the date/text styles match the existing "history" rows exactly.
$ws2.Range("A25:C25").Copy() | Out-Null
$ws2.Range("A26:C26").PasteSpecial(-4122) | Out-Null

# New history entry: 2/13/2017, "Added PDMS", "JEH"
# (Setting the text value here first means "Added PDMS" becomes the next
# unique shared string before any of the PDMS oligomer names below.)
$d = Get-Date -Year 2017 -Month 2 -Day 13 -Hour 0 -Minute 0 -Second 0
$ws2.Range("A26").Value = $d
$ws2.Range("B26").Value = "Added PDMS"
$ws2.Range("C26").Value = "JEH"

# --- LOBSTAHS_rt.windows sheet: append the PDMS oligomer contaminant series
# (PDMS6 .. PDMS27), each with rt_win_max = 30 and rt_win_min = 5, right
# aligned like the other "basic component" rows above them.
For ($i = 0; $i -le 21; $i++) {
  $row = 73 + $i
  $pdmsNum = 6 + $i
  $cellA = $ws1.Range("A" + $row)
  $cellA.Value = "PDMS" + $pdmsNum
  $cellA.HorizontalAlignment = -4152
  $ws1.Range("B" + $row).Value = 30
  $ws1.Range("C" + $row).Value = 5
}

# --- Final view state: LOBSTAHS_rt.windows becomes the active sheet/tab,
# with the last selection left on F82 (Notes tab keeps its own last
# selection at B27, which is set after the data edits above).
$ws2.Range("B27").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("F82").Select() | Out-Null
